$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns so that
# numeric-looking / percentage strings are not auto-converted to numbers.
$rngD = $ws.Range("D2:D51")
$rngE = $ws.Range("E2:E51")
$rngD.NumberFormat = "@"
$rngE.NumberFormat = "@"

$ws.Range("D2").Value = "96.704.81"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.339.44"
$ws.Range("E3").Value = "  -2.38%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "250.06"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "655.31"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "1.40"
$ws.Range("E7").Value = "  -5.16%  "
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "1.01"
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("D11").Value = "3.337.68"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").Value = "40.70"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").Value = "96.476.92"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "6.09"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "0.0000252"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "3.963.28"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "8.71"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "3.332.92"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").Value = "0.559"
$ws.Range("E20").Value = "  +12.48%  "
$ws.Range("D21").Value = "17.43"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "508.15"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("D26").Value = "6.64"
$ws.Range("E26").Value = "  +7.85%  "
$ws.Range("D27").Value = "96.57"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").Value = "12.13"
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("E29").Value = "  -4.40%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "11.33"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").Value = "2.52"
$ws.Range("E33").Value = "  +11.99%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("D36").Value = "28.37"
$ws.Range("E36").Value = "  -4.75%  "
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("D38").Value = "7.80"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").Value = "506.57"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "0.0433"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "0.835"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("D45").Value = "3.68"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "1.68"
$ws.Range("E46").Value = "  +6.91%  "
$ws.Range("D47").Value = "5.58"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "8.50"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "53.72"
$ws.Range("E49").Value = "  +4.83%  "
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").Value = "162.10"
$ws.Range("E51").Value = "  +0.69%  "

# Restore default (unstyled) appearance now that values are committed as text
$rngD.Style = "Normal"
$rngE.Style = "Normal"

